$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new boolean column G: TRUE for row 1, FALSE for rows 2-7
$ws.Range("G1").Value = $true
$ws.Range("G2").Value = $false
$ws.Range("G3").Value = $false
$ws.Range("G4").Value = $false
$ws.Range("G5").Value = $false
$ws.Range("G6").Value = $false
$ws.Range("G7").Value = $false

# Update the active selection to match the diff
$ws.Range("J5").Select()
